$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New ascending period labels (Periodo Mora) for rows 16-31, replacing the
# previous descending 1904 -> 1801 ordering with an ascending 1801 -> 1904
# ordering (database update / "parte 1" of the new account statements).
$periods = @("1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904")

# Corresponding "Valor Mora" (F) amounts, in the same new row order.
$valores = @(29509,29509,29509,29509,29509,29509,29509,29509,31249,31249,31249,31249,31249,31249,31249,16666)

# "Salario Basico" (G) is now a flat updated value for every row.
$salario = 781242

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
    $ws.Range("G$row").Value = $salario
}
